{"js": "// Update the \"Nomination submission\" paragraph: drop the \"TBD\" submission\n// deadline sentence, leaving the subject-line sentence ending in a period\n// right after \"AISEC2025: Best PhD Dissertation Award Nomination\".\n//\n// Before: \"...AISEC2025: Best PhD Dissertation Award Nomination The submission\n//          deadline TBD (all time zones).   \"\n// After:  \"...AISEC2025: Best PhD Dissertation Award Nomination.   \"\n//\n// The trailing \".\" already exists in the document (it ends \"(all time\n// zones).\"), so we only need to delete the text in between \u2014 this keeps the\n// existing plain-formatted \".\" run intact and leaves the bold\n// \"AISEC2025: ... Nomination\" run's formatting untouched.\n\nconst body = context.document.body;\n\nconst target = body.search(\n  \" The submission deadline TBD (all time zones)\",\n  { matchCase: true, matchWholeWord: false }\n);\ntarget.load(\"items/text\");\nawait context.sync();\n\nif (target.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for the deadline sentence, found \" +\n      target.items.length\n  );\n}\n\ntarget.items[0].insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the \"Nomination submission\" paragraph: drop the \"TBD\" submission\n# deadline sentence, leaving the subject-line sentence ending in a period\n# right after \"AISEC2025: Best PhD Dissertation Award Nomination\".\n#\n# Before: \"...AISEC2025: Best PhD Dissertation Award Nomination The submission\n#          deadline TBD (all time zones).   \"\n# After:  \"...AISEC2025: Best PhD Dissertation Award Nomination.   \"\n#\n# The trailing \".\" already exists in the document (it ends \"(all time\n# zones).\"), so we only need to delete the text in between - this keeps the\n# existing plain-formatted \".\" run intact and leaves the bold\n# \"AISEC2025: ... Nomination\" run's formatting untouched.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\" The submission deadline TBD (all time zones)\", $true)\n\nif (-not $found) {\n    throw \"Could not find the deadline sentence to remove\"\n}\n\n$rng.Text = \"\"\n"}
